# Daily update at 8 AM UTC
# Append the next day's row (row 73) to the Wins Over Time data and move the
# "latest row" date-only number format down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day's data (row 73)
$ws.Range("A73").Value = 45813
$ws.Range("B73").Value = 311
$ws.Range("C73").Value = 308
$ws.Range("D73").Value = 314

# Row 72's date cell reverts to the standard datetime format used by the rest
# of the column; the new last row (73) takes on the "latest row" date-only
# format that row 72 used to have.
$ws.Range("A72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A73").NumberFormat = "YYYY-MM-DD"
